$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.257.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.63%  "
$ws.Range("D3").Value = "'2.796.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'116.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("D6").Value = "'341.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("D7").Value = "'0.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.55%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.20%  "
$ws.Range("D10").Value = "'42.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.59%  "
$ws.Range("D11").Value = "'0.0868"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.04%  "
$ws.Range("D12").Value = "'20.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "'3.233.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.99%  "
$ws.Range("D16").Value = "'2.772.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.97%  "
$ws.Range("D17").Value = "'0.888"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("D18").Value = "'52.107.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.30%  "
$ws.Range("D19").Value = "'3.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.83%  "
$ws.Range("D20").Value = "'13.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "'278.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.00%  "
$ws.Range("D24").Value = "'70.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "'2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.07%  "
$ws.Range("D26").Value = "'26.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.42%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'10.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "'34.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'50.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("E33").Value = "  +5.37%  "
$ws.Range("D34").Value = "'0.0829"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "'19.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").Value = "'3.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.22%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0379"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.73%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.89%  "
$ws.Range("D42").Value = "'23.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("E43").Value = "  +4.54%  "
$ws.Range("D44").Value = "'2.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").Value = "'124.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("D46").Value = "'2.104.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "'3.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("E49").Value = "  +7.40%  "
$ws.Range("D50").Value = "'0.909"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +21.90%  "
$ws.Range("D51").Value = "'9.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.76%  "
